$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 207, shifting existing rows 207..239 down to 208..240
$ws.Rows.Item(207).Insert()

# Populate the newly inserted row 207 with the new data record
$ws.Range("A207").Value = 3
$ws.Range("B207").Value = "Femacal de La Calera"
$ws.Range("C207").Value = "Coquimbo"
$ws.Range("D207").Value = 44505
$ws.Range("E207").Value = 5
$ws.Range("F207").Value = 100112031
$ws.Range("G207").Value = "Poroto verde"
$ws.Range("H207").Value = "Magnum"
$ws.Range("I207").Value = "Primera"
$ws.Range("J207").Value = 65
$ws.Range("K207").Value = 41000
$ws.Range("L207").Value = 42000
$ws.Range("M207").Value = 41538
$ws.Range("N207").Value = "`$/malla 25 kilos"
$ws.Range("O207").Value = "Provincia de Limarí"
$ws.Range("P207").Value = 1662
$ws.Range("Q207").Value = 25
$ws.Range("R207").Value = "Hortaliza"
